# Weekly update: a new price observation (week of 2023-01-16) is inserted
# as a new data row right above the existing row 10, pushing the rows that
# used to be 10-14 down to 11-15 (same as Excel's native "Insert Row").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 10; this shifts the old
# rows 10..14 down to 11..15 and keeps their contents/formatting intact.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new observation.
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 44942
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100101
$ws.Cells.Item(10, 8).Value = "Berries"
$ws.Cells.Item(10, 9).Value = 100101001
$ws.Cells.Item(10, 10).Value = "Arándano (blue)"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 60
$ws.Cells.Item(10, 14).Value = 2500
$ws.Cells.Item(10, 15).Value = 2500
$ws.Cells.Item(10, 16).Value = 2500
$ws.Cells.Item(10, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(10, 19).Value = 1250
$ws.Cells.Item(10, 20).Value = 2
